$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the title text typo (B2, part of merged B2:G3) ---
$ws.Range("B2").Value = "This file will provide an overview about the properties of the materials that will be added by the mod"

# --- Row 7: headers for the two tables ---
$ws.Range("B7").Value = "BRUTE  RESOURCES"
$ws.Range("C7").Value = "CLASS"
$ws.Range("D7").Value = "VALUE"
$ws.Range("E7").Value = "RARITY"
$ws.Range("F7").Value = "WEIGHT"

$ws.Range("H7").Value = "REFINED MATERIAL"
$ws.Range("I7").Value = "REQUIRED RESOURCES"
$ws.Range("J7").Value = "VALUE"
$ws.Range("K7").Value = "RARITY"
$ws.Range("L7").Value = "WEIGHT"

# --- Seed rows 8-19 already in "sorted by Class,Name" order so that applying
#     a real Sort (needed to materialise the <sortState> metadata Excel
#     writes out) is a data no-op; we then overwrite with the final curated
#     (non-alphabetical) order that the workbook actually ships with. ---
$seedNames  = @("Dragon Blood","Raw Ancestor Moth Silk","Raw Spidersilk","Aetherium Shard","Adamantium Ore","Ferrous Salts","Manganese Ore","Nickel Ore","Starmetal ore","Golden Sap","Purple Sap","Paper")
$seedClasses = @("Animal Material","Animal Product","Animal Product","Crystaline Mineral","Metallic Mineral","Metallic Mineral","Metallic Mineral","Metallic Mineral","Metallic Mineral","Plant Material","Plant Material","Produce")

$r = 8
for ($i = 0; $i -lt $seedNames.Count; $i++) {
  $ws.Cells.Item($r, 2).Value = $seedNames[$i]
  $ws.Cells.Item($r, 3).Value = $seedClasses[$i]
  $r = $r + 1
}

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("C8:C19"))
$sortObj.SortFields.Add($ws.Range("B8:B19"))
$sortObj.SetRange($ws.Range("B8:C19"))
$sortObj.Header = -4142
$sortObj.Apply()

# --- Rows 8-19: final curated brute-resources/class listing ---
$ws.Range("B8").Value = "Adamantium Ore"
$ws.Range("C8").Value = "Metallic Mineral"

$ws.Range("B9").Value = "Aetherium Shard"
$ws.Range("C9").Value = "Crystaline Mineral"

$ws.Range("B10").Value = "Ferrous Salts"
$ws.Range("C10").Value = "Metallic Mineral"

$ws.Range("B11").Value = "Manganese Ore"
$ws.Range("C11").Value = "Metallic Mineral"

$ws.Range("B12").Value = "Nickel Ore"
$ws.Range("C12").Value = "Metallic Mineral"

$ws.Range("B13").Value = "Starmetal ore"
$ws.Range("C13").Value = "Metallic Mineral"

$ws.Range("B14").Value = "Dragon Blood"
$ws.Range("C14").Value = "Animal Material"

$ws.Range("B15").Value = "Golden Sap"
$ws.Range("C15").Value = "Plant Material"

$ws.Range("B16").Value = "Purple Sap"
$ws.Range("C16").Value = "Plant Material"

$ws.Range("B17").Value = "Raw Ancestor Moth Silk"
$ws.Range("C17").Value = "Animal Product"

$ws.Range("B18").Value = "Raw Spidersilk"
$ws.Range("C18").Value = "Animal Product"

$ws.Range("B19").Value = "Paper"
$ws.Range("C19").Value = "Produce"

# --- Rows 20-28: additional brute resources added after the sort ---
$ws.Range("B20").Value = "Rough Oak"
$ws.Range("C20").Value = "Plant Material"

$ws.Range("B21").Value = "Rough Nightwood"
$ws.Range("C21").Value = "Plant Material"

$ws.Range("B22").Value = "Rough Maple"
$ws.Range("C22").Value = "Plant Material"

$ws.Range("B23").Value = "Rough Birch"
$ws.Range("C23").Value = "Plant Material"

$ws.Range("B24").Value = "Raw Ebonthread"
$ws.Range("C24").Value = "Plant Material"

$ws.Range("B25").Value = "Raw Akaviri Silk"
$ws.Range("C25").Value = "Animal Product"

$ws.Range("B26").Value = "Raw Cotton"
$ws.Range("C26").Value = "Plant Material"

$ws.Range("B27").Value = "Raw Ironweed"
$ws.Range("C27").Value = "Plant Material"

$ws.Range("B28").Value = "Raw Void Bloom"
$ws.Range("C28").Value = "Plant Material"

# --- Column widths (best-fit like); the engine quantises stored widths to
#     1/6-character increments, so these inputs are chosen to land on the
#     closest achievable stored width to the target (22.140625, 17.42578125,
#     18, 20.7109375). ---
$ws.Columns.Item(2).ColumnWidth = 21.333333333333332
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668
$ws.Columns.Item(8).ColumnWidth = 17.166666666666668
$ws.Columns.Item(9).ColumnWidth = 19.833333333333332

# --- Row 2 explicit custom height (matches customHeight="1" in target) ---
$ws.Rows.Item(2).RowHeight = 15.75

# --- Selection moved to F11 ---
$null = $ws.Range("F11").Select()

Write-Host "done"
